# Add new columns I ("I0") and J ("IF") to the stats sheet, mirroring the
# header style already used by the other header cells (H1, etc.), and fill
# in the per-row numeric values for rows 2..71.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the existing header formatting (bold font, borders, centered
# alignment) from H1 onto the two new header cells so they match the rest
# of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows -----------------------------------------------------------
$rows = @(
    @(2, 9, 9),
    @(3, 9, 9),
    @(4, 7, 8),
    @(5, 8, 9),
    @(6, 9, 9),
    @(7, 8, 8),
    @(8, 9, 9),
    @(9, 9, 9),
    @(10, 8, 9),
    @(11, 9, 9),
    @(12, 9, 9),
    @(13, 9, 9),
    @(14, 10, 10),
    @(15, 9, 9),
    @(16, 9, 9),
    @(17, 9, 9),
    @(18, 9, 9),
    @(19, 10, 10),
    @(20, 10, 10),
    @(21, 9, 9),
    @(22, 9, 9),
    @(23, 9, 9),
    @(24, 8, 8),
    @(25, 8, 8),
    @(26, 8, 8),
    @(27, 8, 8),
    @(28, 10, 10),
    @(29, 8, 8),
    @(30, 8, 8),
    @(31, 8, 8),
    @(32, 8, 8),
    @(33, 8, 8),
    @(34, 8, 8),
    @(35, 8, 8),
    @(36, 8, 8),
    @(37, 8, 8),
    @(38, 8, 8),
    @(39, 9, 9),
    @(40, 8, 8),
    @(41, 8, 8),
    @(42, 8, 8),
    @(43, 8, 8),
    @(44, 8, 8),
    @(45, 9, 9),
    @(46, 8, 8),
    @(47, 9, 9),
    @(48, 9, 9),
    @(49, 9, 9),
    @(50, 8, 8),
    @(51, 9, 9),
    @(52, 8, 8),
    @(53, 8, 8),
    @(54, 8, 8),
    @(55, 8, 8),
    @(56, 9, 9),
    @(57, 9, 9),
    @(58, 8, 8),
    @(59, 8, 8),
    @(60, 8, 8),
    @(61, 8, 8),
    @(62, 8, 8),
    @(63, 8, 8),
    @(64, 8, 8),
    @(65, 8, 8),
    @(66, 8, 8),
    @(67, 8, 8),
    @(68, 7, 7),
    @(69, 7, 7),
    @(70, 5, 5),
    @(71, 3, 3)
)

foreach ($r in $rows) {
    $row = $r[0]
    $iVal = $r[1]
    $jVal = $r[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
